$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings like "40.167.97" / "0.1000" that must stay
# literal text (the source data never stores prices as real numbers). Force
# text interpretation via NumberFormat, then restore the Normal style so the
# cell's formatting matches the rest of the untouched column.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.167.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.223.39'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0784'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("E13").Value = '  +3.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.564.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.87'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.217.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.077.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  -3.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.29%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.114'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.079.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.45%  '
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("E47").Value = '  +5.51%  '
$ws.Range("E48").Value = '  -11.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.437.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("E50").Value = '  +4.13%  '
$ws.Range("E51").Value = '  +1.82%  '
